$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the used range contents (but keep formatting) so the shared-string pool
# is rebuilt fresh in the exact column-major order the target file expects.
$ws.Range("A1:G21").ClearContents()

# Force text storage (not auto-converted numbers) for the rank column and the
# numeric-looking stat columns, matching the original file where every cell is
# a shared string, not a number. (Applied as two single-area ranges -- a
# multi-area "A2:A21,C2:G21" union does not reliably propagate NumberFormat
# to every area.)
$rankRange = $ws.Range("A2:A21")
$statsRange = $ws.Range("C2:G21")
$rankRange.NumberFormat = "@"
$statsRange.NumberFormat = "@"

# --- Header row (row 1) ---
$ws.Range("A1").Value = "#"
$ws.Range("B1").Value = "Equipe"
$ws.Range("C1").Value = "Cartões"
$ws.Range("D1").Value = "Escanteios"
$ws.Range("E1").Value = "1.5+"
$ws.Range("F1").Value = "2.5+"
$ws.Range("G1").Value = "Med. Gols"

# --- Data written column-by-column (A, then B, then C..G) ---
# Column A
$ws.Range("A2").Value = "1."
$ws.Range("A3").Value = "2."
$ws.Range("A4").Value = "3."
$ws.Range("A5").Value = "4."
$ws.Range("A6").Value = "5."
$ws.Range("A7").Value = "6."
$ws.Range("A8").Value = "7."
$ws.Range("A9").Value = "8."
$ws.Range("A10").Value = "9."
$ws.Range("A11").Value = "10."
$ws.Range("A12").Value = "11."
$ws.Range("A13").Value = "12."
$ws.Range("A14").Value = "13."
$ws.Range("A15").Value = "14."
$ws.Range("A16").Value = "15."
$ws.Range("A17").Value = "16."
$ws.Range("A18").Value = "17."
$ws.Range("A19").Value = "18."
$ws.Range("A20").Value = "19."
$ws.Range("A21").Value = "20."

# Column B
$ws.Range("B2").Value = "Arsenal"
$ws.Range("B3").Value = "Liverpool"
$ws.Range("B4").Value = "Man City"
$ws.Range("B5").Value = "Aston Villa"
$ws.Range("B6").Value = "Tottenham"
$ws.Range("B7").Value = "Newcastle"
$ws.Range("B8").Value = "Manchester United"
$ws.Range("B9").Value = "West Ham"
$ws.Range("B10").Value = "Chelsea"
$ws.Range("B11").Value = "Brighton"
$ws.Range("B12").Value = "Wolves"
$ws.Range("B13").Value = "Fulham"
$ws.Range("B14").Value = "Bournemouth"
$ws.Range("B15").Value = "Crystal Palace"
$ws.Range("B16").Value = "Brentford"
$ws.Range("B17").Value = "Everton"
$ws.Range("B18").Value = "Nottingham"
$ws.Range("B19").Value = "Luton"
$ws.Range("B20").Value = "Burnley"
$ws.Range("B21").Value = "Sheffield Utd"

# Column C
$ws.Range("C2").Value = "1.6"
$ws.Range("C3").Value = "1.8"
$ws.Range("C4").Value = "1.6"
$ws.Range("C5").Value = "2.5"
$ws.Range("C6").Value = "2.5"
$ws.Range("C7").Value = "2.0"
$ws.Range("C8").Value = "2.3"
$ws.Range("C9").Value = "2.3"
$ws.Range("C10").Value = "3.0"
$ws.Range("C11").Value = "2.4"
$ws.Range("C12").Value = "2.6"
$ws.Range("C13").Value = "2.2"
$ws.Range("C14").Value = "2.1"
$ws.Range("C15").Value = "1.9"
$ws.Range("C16").Value = "2.3"
$ws.Range("C17").Value = "2.2"
$ws.Range("C18").Value = "2.3"
$ws.Range("C19").Value = "1.8"
$ws.Range("C20").Value = "1.9"
$ws.Range("C21").Value = "2.8"

# Column D
$ws.Range("D2").Value = "7.0"
$ws.Range("D3").Value = "7.4"
$ws.Range("D4").Value = "7.8"
$ws.Range("D5").Value = "6.4"
$ws.Range("D6").Value = "6.0"
$ws.Range("D7").Value = "4.9"
$ws.Range("D8").Value = "5.9"
$ws.Range("D9").Value = "4.3"
$ws.Range("D10").Value = "5.3"
$ws.Range("D11").Value = "5.6"
$ws.Range("D12").Value = "4.1"
$ws.Range("D13").Value = "5.8"
$ws.Range("D14").Value = "6.2"
$ws.Range("D15").Value = "4.6"
$ws.Range("D16").Value = "4.7"
$ws.Range("D17").Value = "4.8"
$ws.Range("D18").Value = "3.9"
$ws.Range("D19").Value = "5.4"
$ws.Range("D20").Value = "4.8"
$ws.Range("D21").Value = "3.6"

# Column E
$ws.Range("E2").Value = "79%"
$ws.Range("E3").Value = "91%"
$ws.Range("E4").Value = "81%"
$ws.Range("E5").Value = "88%"
$ws.Range("E6").Value = "97%"
$ws.Range("E7").Value = "81%"
$ws.Range("E8").Value = "78%"
$ws.Range("E9").Value = "94%"
$ws.Range("E10").Value = "88%"
$ws.Range("E11").Value = "87%"
$ws.Range("E12").Value = "85%"
$ws.Range("E13").Value = "79%"
$ws.Range("E14").Value = "91%"
$ws.Range("E15").Value = "79%"
$ws.Range("E16").Value = "85%"
$ws.Range("E17").Value = "76%"
$ws.Range("E18").Value = "85%"
$ws.Range("E19").Value = "88%"
$ws.Range("E20").Value = "91%"
$ws.Range("E21").Value = "91%"

# Column F
$ws.Range("F2").Value = "62%"
$ws.Range("F3").Value = "67%"
$ws.Range("F4").Value = "63%"
$ws.Range("F5").Value = "71%"
$ws.Range("F6").Value = "81%"
$ws.Range("F7").Value = "72%"
$ws.Range("F8").Value = "66%"
$ws.Range("F9").Value = "65%"
$ws.Range("F10").Value = "69%"
$ws.Range("F11").Value = "62%"
$ws.Range("F12").Value = "64%"
$ws.Range("F13").Value = "65%"
$ws.Range("F14").Value = "70%"
$ws.Range("F15").Value = "55%"
$ws.Range("F16").Value = "68%"
$ws.Range("F17").Value = "45%"
$ws.Range("F18").Value = "53%"
$ws.Range("F19").Value = "71%"
$ws.Range("F20").Value = "59%"
$ws.Range("F21").Value = "69%"

# Column G
$ws.Range("G2").Value = "3.18"
$ws.Range("G3").Value = "3.24"
$ws.Range("G4").Value = "3.38"
$ws.Range("G5").Value = "3.56"
$ws.Range("G6").Value = "3.56"
$ws.Range("G7").Value = "3.78"
$ws.Range("G8").Value = "2.97"
$ws.Range("G9").Value = "3.44"
$ws.Range("G10").Value = "3.69"
$ws.Range("G11").Value = "3.19"
$ws.Range("G12").Value = "3.00"
$ws.Range("G13").Value = "3.06"
$ws.Range("G14").Value = "3.27"
$ws.Range("G15").Value = "2.97"
$ws.Range("G16").Value = "3.26"
$ws.Range("G17").Value = "2.48"
$ws.Range("G18").Value = "3.00"
$ws.Range("G19").Value = "3.59"
$ws.Range("G20").Value = "3.12"
$ws.Range("G21").Value = "3.61"

# Drop the temporary text-format styling so the saved styles.xml / cell "s"
# attributes match the original (un-styled) data cells.
$rankRange.ClearFormats()
$statsRange.ClearFormats()